$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: new day entry (2020.05.11 / "Priority is to sort out git submodules!") ---

# A25 needs to stay plain text "2020.05.11" (like the other date-label cells
# in column A) rather than being auto-parsed into a serial date. Briefly
# force a Text number format while the value is typed in, then restore the
# default ("Normal") cell style so no stray custom format lingers on it.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "2020.05.11"
$ws.Range("A25").Style = "Normal"

# B25 is a clock-in time, formatted the same way as the other B-column cells.
$ws.Range("B25").NumberFormat = $ws.Range("B24").NumberFormat
$ws.Range("B25").Value = 0.49305555555555558

# E25 holds the day's note, wrapped like the rest of column E.
$ws.Range("E25").WrapText = $true
$ws.Range("E25").Value = "Priority is to sort out git submodules!"

# D25 already carries the shared "C-B" formula (shared formula si=0,
# range D7:D26) and recalculates automatically once B25/C25 change.

# --- Selection / scroll position bookkeeping (matches the saved view state) ---
$ws.Range("A26").Select()
